# Add a new row (row 7) of timing data to the "tijden" sheet, matching the
# author's commit: "indentfoutje eruit gehaald + excel sheet geupdatet"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row ------------------------------------------------------
$ws.Range("A7").Value = "Bij move: tuple->list->bewerken->tuple (dus maar 1x naar list en terug)"
$ws.Range("B7").Value = 0.002
$ws.Range("C7").Value = 3.35
$ws.Range("D7").Value = 1
$ws.Range("D7").NumberFormat = "0.00"
$ws.Range("E7").Value = 0.19

# --- Widen column A slightly so the longer text fits -------------------
$ws.Columns.Item(1).ColumnWidth = 60.66796875

# --- Selection moves to B7 ----------------------------------------------
[void]$ws.Range("B7").Select()
